# Add season record columns (Wins, Losses, Ties) to the KCR_1993 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row (data rows start at row 2).
$lastRow = $ws.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count

# Header cells (AD1:AF1) - styled like the other header cells (bold, centered, bordered).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from the neighboring header cell (AC1) onto the new headers
# so they match the existing bold / centered / bordered header style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Season record values for every player row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 78
    $ws.Cells.Item($r, 32).Value = 0
}
